# Updates cryptos list data (prices / 1h volume changes) per upstream scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as a plain number must be forced back to Text
# so they keep matching the "Price"/"Volume(1h)" columns existing string format
# (avoids silently turning "321.55" into a numeric 321.55 cell).
$textFormatCells = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D12",
    "D14",
    "D16",
    "D19",
    "D22",
    "D23",
    "D24",
    "D25",
    "D27",
    "D28",
    "D30",
    "D31",
    "D32",
    "D33",
    "D39",
    "D40",
    "D41",
    "D42",
    "D45",
    "D46",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the refreshed values scraped by the GitHub Actions job.
$ws.Range('D2').Value = '44.447.52'
$ws.Range('E2').Value = '  +3.86%  '
$ws.Range('D3').Value = '2.274.38'
$ws.Range('E3').Value = '  +3.20%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '321.55'
$ws.Range('E5').Value = '  +2.04%  '
$ws.Range('D6').Value = '105.25'
$ws.Range('E6').Value = '  +6.57%  '
$ws.Range('D7').Value = '0.591'
$ws.Range('E7').Value = '  +0.55%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = '0.572'
$ws.Range('E9').Value = '  +2.39%  '
$ws.Range('D10').Value = '38.79'
$ws.Range('E10').Value = '  +5.99%  '
$ws.Range('E11').Value = '  +2.31%  '
$ws.Range('D12').Value = '7.88'
$ws.Range('E12').Value = '  +2.87%  '
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('D14').Value = '0.885'
$ws.Range('E14').Value = '  +3.56%  '
$ws.Range('D15').Value = '2.617.83'
$ws.Range('E15').Value = '  +3.21%  '
$ws.Range('D16').Value = '14.58'
$ws.Range('E16').Value = '  +2.68%  '
$ws.Range('D17').Value = '2.273.45'
$ws.Range('E17').Value = '  +3.01%  '
$ws.Range('D18').Value = '44.325.40'
$ws.Range('E18').Value = '  +3.89%  '
$ws.Range('D19').Value = '13.91'
$ws.Range('E19').Value = '  -3.65%  '
$ws.Range('E20').Value = '  +4.69%  '
$ws.Range('E21').Value = '  +2.17%  '
$ws.Range('D22').Value = '66.51'
$ws.Range('E22').Value = '  +2.23%  '
$ws.Range('D23').Value = '3.21'
$ws.Range('E23').Value = '  +2.08%  '
$ws.Range('D24').Value = '239.98'
$ws.Range('E24').Value = '  +1.87%  '
$ws.Range('D25').Value = '2.23'
$ws.Range('E25').Value = '  +5.58%  '
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('D27').Value = '10.23'
$ws.Range('E27').Value = '  +2.71%  '
$ws.Range('D28').Value = '38.57'
$ws.Range('E28').Value = '  +13.28%  '
$ws.Range('E29').Value = '  -0.63%  '
$ws.Range('D30').Value = '6.52'
$ws.Range('E30').Value = '  +3.79%  '
$ws.Range('D31').Value = '20.70'
$ws.Range('E31').Value = '  +1.13%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').Value = '162.39'
$ws.Range('E32').Value = '  +4.99%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.0887'
$ws.Range('E33').Value = '  +0.71%  '
$ws.Range('E34').Value = '  -0.78%  '
$ws.Range('E35').Value = '  +8.01%  '
$ws.Range('E36').Value = '  +6.28%  '
$ws.Range('E37').Value = '  +1.96%  '
$ws.Range('E38').Value = '  +0.56%  '
$ws.Range('D39').Value = '3.95'
$ws.Range('E39').Value = '  +4.14%  '
$ws.Range('D40').Value = '4.45'
$ws.Range('E40').Value = '  +1.06%  '
$ws.Range('D41').Value = '15.76'
$ws.Range('E41').Value = '  +28.30%  '
$ws.Range('D42').Value = '0.0329'
$ws.Range('E42').Value = '  +1.77%  '
$ws.Range('E43').Value = '  +0.32%  '
$ws.Range('D44').Value = '1.781.96'
$ws.Range('E44').Value = '  -4.86%  '
$ws.Range('D45').Value = '0.210'
$ws.Range('E45').Value = '  +1.65%  '
$ws.Range('D46').Value = '86.35'
$ws.Range('E46').Value = '  -1.40%  '
$ws.Range('E47').Value = '  +2.25%  '
$ws.Range('D48').Value = '60.87'
$ws.Range('E48').Value = '  +0.46%  '
$ws.Range('D49').Value = '75.42'
$ws.Range('E49').Value = '  +0.24%  '
$ws.Range('D50').Value = '1.72'
$ws.Range('E50').Value = '  +9.05%  '
$ws.Range('D51').Value = '104.34'
$ws.Range('E51').Value = '  +2.17%  '
